# SCD0018-004 - CRO atau CRM mengajukan data Non Sales.xlsx
# Update TC_ID from "DGS-296" to "SCD0018-004" and rename the sheet/tab
# to match the new TC_ID (SCD0018). Also widen the TC_ID column so the
# longer id fits, and move the active cell selection to B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet tab: SCD0281 -> SCD0018
$ws.Name = "SCD0018"

# Update TC_ID column (B) for all 3 test-case rows
$ws.Range("B2").Value = "SCD0018-004"
$ws.Range("B3").Value = "SCD0018-004"
$ws.Range("B4").Value = "SCD0018-004"

# Widen column B to fit the new, longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.7

# Move the current selection to B5
$ws.Range("B5").Select()
